{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 2 (\"For our first meeting, ...coming together. \") gets new\n// sentences appended about Tim, Hugo and Taylen's contributions.\nconst p2 = paragraphs.items[1];\np2.insertText(\n  \"Tim and Hugo have done an exceptional job with the project artifact. \" +\n  \"From assessment two, they have taken the project idea and brought it to life \" +\n  \"as something that can be marketable in the real world. Taylen has done an \" +\n  \"amazing job coordinating the team and have all worked together to come up \" +\n  \"with an amazing piece of work.\",\n  Word.InsertLocation.end\n);\n\n// Paragraph 3 (\"Through my group work ...\") becomes \"Throughout my group work ...\"\nconst p3 = paragraphs.items[2];\nconst searchResults = p3.search(\"Through\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"text\");\nawait context.sync();\n\nsearchResults.items[0].insertText(\"out\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 2 (\"For our first meeting, ...coming together. \") gets new\n# sentences appended about Tim, Hugo and Taylen's contributions.\n$p2 = $d.Paragraphs.Item(2)\n$r2 = $p2.Range\n$r2.MoveEnd(1, -1) | Out-Null\n$r2.Collapse(0) | Out-Null\n$r2.InsertAfter(\"Tim and Hugo have done an exceptional job with the project artifact. From assessment two, they have taken the project idea and brought it to life as something that can be marketable in the real world. Taylen has done an amazing job coordinating the team and have all worked together to come up with an amazing piece of work.\")\n\n# Paragraph 3 (\"Through my group work ...\") becomes \"Throughout my group work ...\"\n$p3 = $d.Paragraphs.Item(3)\n$r3 = $p3.Range\n$r3.Find.Execute(\"Through\", $false, $true, $false, $false, $false, $true, 1, $false, \"Throughout\", 1) | Out-Null\n"}
